$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells retain text formatting (values are stored as text strings, e.g. "277.90", "3.85%")
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '277.90'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.85%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.80'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.32%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.917'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.77%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06388'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.93%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.992'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.78%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.353'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.76%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.352'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '49.04%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8856'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.19%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1481'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.83%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05193'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.59%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07422'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '4.63%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03150'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.40%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09065'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.61%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001561'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.06%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006338'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '4.66%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006013'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.65%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.86%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.19%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1332'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.51%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.901'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.99%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04357'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.92%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001177'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.34%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003677'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-11.07%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001201'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.06%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001618'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-3.75%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04073'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.09%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006646'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '58.24%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1175'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '5.59%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002362'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '11.38%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01282'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.78%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005266'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.75%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.02%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '1,788.43%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02123'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-13.25%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.02%'
